$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I4").Value = 2
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48
